# Updates cryptocurrency price/volume figures scraped on Fri Feb  2 23:38:47 UTC 2024.
# Two rows (Monero / InjectiveProtocol) also swapped position-for-position in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.167.53"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").Value = "2.306.71"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("E5").Value = "  +0.54%  "

# Row 6
$ws.Range("D6").Value = "'100.61"
$ws.Range("E6").Value = "  +3.17%  "

# Row 7
$ws.Range("E7").Value = "  +0.66%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +1.69%  "

# Row 10
$ws.Range("D10").Value = "'36.71"
$ws.Range("E10").Value = "  +8.63%  "

# Row 11
$ws.Range("E11").Value = "  +0.06%  "

# Row 12
$ws.Range("E12").Value = "  +0.73%  "

# Row 13
$ws.Range("D13").Value = "'17.70"
$ws.Range("E13").Value = "  +3.55%  "

# Row 14
$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  +2.17%  "

# Row 15
$ws.Range("D15").Value = "2.665.38"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").Value = "2.306.47"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("E17").Value = "  -1.54%  "

# Row 18
$ws.Range("D18").Value = "43.070.55"
$ws.Range("E18").Value = "  +0.30%  "

# Row 19
$ws.Range("D19").Value = "'12.73"
$ws.Range("E19").Value = "  +9.40%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  +0.26%  "

# Row 21
$ws.Range("E21").Value = "  +1.76%  "

# Row 22
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").Value = "'236.23"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +7.33%  "

# Row 25
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  -0.12%  "

# Row 26
$ws.Range("D26").Value = "'0.993"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("E27").Value = "  +3.21%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'34.61"
$ws.Range("E28").Value = "  +2.25%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'167.80"
$ws.Range("E29").Value = "  +0.65%  "

# Row 30
$ws.Range("E30").Value = "  -1.03%  "

# Row 31
$ws.Range("E31").Value = "  +0.34%  "

# Row 32
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("D33").Value = "'5.05"
$ws.Range("E33").Value = "  +1.74%  "

# Row 34
$ws.Range("D34").Value = "'17.64"
$ws.Range("E34").Value = "  +4.71%  "

# Row 35
$ws.Range("D35").Value = "'4.60"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36
$ws.Range("E36").Value = "  -0.88%  "

# Row 37
$ws.Range("D37").Value = "'0.0691"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38
$ws.Range("E38").Value = "  +1.69%  "

# Row 39
$ws.Range("E39").Value = "  +0.81%  "

# Row 40
$ws.Range("E40").Value = "  -0.27%  "

# Row 41
$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = "  +0.52%  "

# Row 42
$ws.Range("E42").Value = "  +3.22%  "

# Row 43
$ws.Range("D43").Value = "1.982.93"
$ws.Range("E43").Value = "  -0.48%  "

# Row 44
$ws.Range("E44").Value = "  -4.11%  "

# Row 45
$ws.Range("D45").Value = "'10.26"
$ws.Range("E45").Value = "  +4.54%  "

# Row 46
$ws.Range("D46").Value = "'17.82"
$ws.Range("E46").Value = "  +1.38%  "

# Row 47
$ws.Range("E47").Value = "  +2.10%  "

# Row 48
$ws.Range("D48").Value = "'55.31"
$ws.Range("E48").Value = "  +3.68%  "

# Row 49
$ws.Range("D49").Value = "'1.55"
$ws.Range("E49").Value = "  +4.10%  "

# Row 50
$ws.Range("E50").Value = "  +0.32%  "

# Row 51
$ws.Range("D51").Value = "'70.82"
$ws.Range("E51").Value = "  +1.19%  "
